# The commit swaps the deck's theme color palette from the "Integral" /
# "Red Violet" scheme over to the stock PowerPoint "Office Theme" / "Office"
# color scheme (what theme1.xml looks like in the target diff).
#
# PowerPoint's automation object model exposes the live theme palette via
# SlideMaster.Theme.ThemeColorScheme.Colors(n).RGB (n = 1..12, in the same
# order the colors appear in a:clrScheme: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). Setting each entry rewrites the corresponding
# <a:srgbClr val="…"/> in the theme part used by the slide master.
# (.RGB uses the standard OLE/VBA packing: R + G*256 + B*65536.)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$cs = $theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0x00 + (0x00 * 256) + (0x00 * 65536)   # dk1      000000
$cs.Colors(2).RGB  = 0xFF + (0xFF * 256) + (0xFF * 65536)   # lt1      FFFFFF
$cs.Colors(3).RGB  = 0x44 + (0x54 * 256) + (0x6A * 65536)   # dk2      44546A
$cs.Colors(4).RGB  = 0xE7 + (0xE6 * 256) + (0xE6 * 65536)   # lt2      E7E6E6
$cs.Colors(5).RGB  = 0x5B + (0x9B * 256) + (0xD5 * 65536)   # accent1  5B9BD5
$cs.Colors(6).RGB  = 0xED + (0x7D * 256) + (0x31 * 65536)   # accent2  ED7D31
$cs.Colors(7).RGB  = 0xA5 + (0xA5 * 256) + (0xA5 * 65536)   # accent3  A5A5A5
$cs.Colors(8).RGB  = 0xFF + (0xC0 * 256) + (0x00 * 65536)   # accent4  FFC000
$cs.Colors(9).RGB  = 0x44 + (0x72 * 256) + (0xC4 * 65536)   # accent5  4472C4
$cs.Colors(10).RGB = 0x70 + (0xAD * 256) + (0x47 * 65536)   # accent6  70AD47
$cs.Colors(11).RGB = 0x05 + (0x63 * 256) + (0xC1 * 65536)   # hlink    0563C1
$cs.Colors(12).RGB = 0x95 + (0x4F * 256) + (0x72 * 65536)   # folHlink 954F72

Write-Output "Theme color scheme updated to Office Theme palette."
